$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '36.621.86'
$ws.Range("E2").Value = '  -0.84%  '

# Row 3
$ws.Range("D3").Value = '2.100.06'
$ws.Range("E3").Value = '  +9.49%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.60%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.661'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.46%  '

# Row 7
$ws.Range("E7").Value = '  +0.11%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.74'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.80%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '60.34'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.09%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.378'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.04%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0746'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.42%  '

# Row 12
$ws.Range("E12").Value = '  +0.46%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.66'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.21%  '

# Row 14
$ws.Range("D14").Value = '2.407.03'
$ws.Range("E14").Value = '  +9.63%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.836'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.08%  '

# Row 16
$ws.Range("D16").Value = '2.101.70'
$ws.Range("E16").Value = '  +9.62%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.13'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.53%  '

# Row 18
$ws.Range("D18").Value = '36.651.33'
$ws.Range("E18").Value = '  -1.69%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.63%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0835'
$ws.Range("E20").Value = '  -3.41%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.37'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.21%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '240.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.79%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.27'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.18%  '

# Row 24
$ws.Range("E24").Value = '  -0.03%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.19%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.02'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.77%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.26'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +13.12%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.76%  '

# Row 29
$ws.Range("E29").Value = '  -10.12%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '28.49'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +53.34%  '

# Row 31
$ws.Range("E31").Value = '  -5.34%  '

# Row 32
$ws.Range("E32").Value = '  -2.27%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0620'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.12%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.45'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +20.91%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.979'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.72%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0896'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.70%  '

# Row 37
$ws.Range("E37").Value = '  -0.02%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.85'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.92%  '

# Row 39
$ws.Range("E39").Value = '  -5.80%  '

# Row 40
$ws.Range("E40").Value = '  -11.07%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0225'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.62%  '

# Row 42
$ws.Range("E42").Value = '  +5.52%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '98.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.42%  '

# Row 44
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.43'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.53%  '

# Row 45
$ws.Range("B45").Value = 'HuobiToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.53%  '

# Row 46
$ws.Range("D46").Value = '1.340.15'
$ws.Range("E46").Value = '  -0.96%  '

# Row 47
$ws.Range("E47").Value = '  +3.82%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.11'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.49%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.88'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.69%  '

# Row 50
$ws.Range("D50").Value = '2.297.08'
$ws.Range("E50").Value = '  +9.86%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.25'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.49%  '
